$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "42.180.11"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "2.268.02"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue "D5" "304.88"
$ws.Range("E5").Value = "  +0.47%  "
Set-TextValue "D6" "96.59"
$ws.Range("E6").Value = "  +4.25%  "
Set-TextValue "D7" "0.530"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.02%  "
Set-TextValue "D9" "0.491"
$ws.Range("E9").Value = "  +1.25%  "
Set-TextValue "D10" "35.55"
$ws.Range("E10").Value = "  +9.61%  "
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("E12").Value = "  -0.25%  "
Set-TextValue "D13" "6.65"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("D14").Value = "2.619.49"
$ws.Range("E14").Value = "  +0.16%  "
Set-TextValue "D15" "14.45"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").Value = "2.273.80"
$ws.Range("E16").Value = "  +0.27%  "
Set-TextValue "D17" "0.794"
$ws.Range("E17").Value = "  +1.65%  "
$ws.Range("D18").Value = "42.090.15"
$ws.Range("E18").Value = "  +0.79%  "
Set-TextValue "D19" "12.46"
$ws.Range("E19").Value = "  -1.76%  "
$ws.Range("D20").Value = "0.0₃0908"
$ws.Range("E20").Value = "  +0.01%  "
Set-TextValue "D21" "6.00"
$ws.Range("E21").Value = "  +1.08%  "
Set-TextValue "D22" "67.85"
$ws.Range("E22").Value = "  +0.36%  "
Set-TextValue "D23" "238.44"
$ws.Range("E23").Value = "  -2.36%  "
Set-TextValue "D24" "2.57"
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("E26").Value = "  -0.08%  "
Set-TextValue "D27" "23.72"
$ws.Range("E27").Value = "  -1.27%  "
Set-TextValue "D28" "37.26"
$ws.Range("E28").Value = "  +6.56%  "
Set-TextValue "D29" "9.54"
$ws.Range("E29").Value = "  -0.51%  "
Set-TextValue "D30" "2.11"
$ws.Range("E30").Value = "  +1.65%  "
Set-TextValue "D31" "159.45"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("E32").Value = "  -0.42%  "
Set-TextValue "D33" "0.999"
$ws.Range("E33").Value = "  +0.01%  "
Set-TextValue "D34" "3.19"
$ws.Range("E34").Value = "  +5.54%  "
$ws.Range("E35").Value = "  -0.65%  "
Set-TextValue "D36" "17.24"
$ws.Range("E36").Value = "  +1.91%  "
$ws.Range("E39").Value = "  +2.15%  "
$ws.Range("E40").Value = "  -1.51%  "
Set-TextValue "D41" "4.08"
$ws.Range("E41").Value = "  +4.02%  "
$ws.Range("E42").Value = "  +13.40%  "
$ws.Range("D43").Value = "1.990.37"
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("E46").Value = "  +1.16%  "
Set-TextValue "D47" "9.93"
$ws.Range("E47").Value = "  -3.91%  "
Set-TextValue "D48" "53.17"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("E49").Value = "  +0.57%  "
Set-TextValue "D50" "72.21"
Set-TextValue "D51" "91.37"
$ws.Range("E51").Value = "  -0.27%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D37" "2.37"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D38" "0.105"
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D44" "19.08"
$ws.Range("E44").Value = "  -4.45%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D45" "0.0284"
$ws.Range("E45").Value = "  +0.69%  "
